$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new record row at row 294 (shifts existing rows 294:380 down to 295:381)
$ws.Rows(294).Insert()

# Populate the newly inserted row with the new data record
$ws.Range("A294").Value = 4
$ws.Range("B294").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C294").Value = "Los Lagos"
$ws.Range("D294").Value = 44985
$ws.Range("E294").Value = 10
$ws.Range("F294").Value = "Fruta"
$ws.Range("G294").Value = 100108
$ws.Range("H294").Value = "Tropicales y subtropicales"
$ws.Range("I294").Value = 100108005
$ws.Range("J294").Value = "Piña"
$ws.Range("K294").Value = "Caramelo"
$ws.Range("L294").Value = "Primera"
$ws.Range("M294").Value = 200
$ws.Range("N294").Value = 25000
$ws.Range("O294").Value = 26000
$ws.Range("P294").Value = 25500
$ws.Range("Q294").Value = "$/caja 12 unidades"
$ws.Range("R294").Value = "Ecuador"
$ws.Range("S294").Value = 2125
$ws.Range("T294").Value = 12
